$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('3364 Assisi Centre Aged Care Rosanna', 29)
    ,@('3376 Royal Freemasons Coppin Centre Melbourne', 19)
    ,@('3622 Olivet Care Aged Care Services Ringwood', 12)
    ,@('3633 Lifeview Emerald Glades Aged Care Emerald', 10)
    ,@('3825 TLC Forest Lodge Residential Aged Care Frankston North', 14)
    ,@('3961 Heritage Care Water Gardens Aged Care Facility Sydenham', 22)
    ,@('4167 Royal Freemasons Centennial Lodge Wantirna South', 22)
    ,@('4282 Villa Maria Catholic Homes (VMCH) Wantirna At-Home Aged Care Stud Road Wantirna', 11)
    ,@('Alamanda K9 College Point Cook', 10)
    ,@('Antonine College Cedar Campus Coburg', 12)
    ,@('Australian Meat Group Abattoir Dandenong South', 10)
    ,@('Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh', 40)
    ,@('Baden Powell College Tarneit', 13)
    ,@('Chairo Christian School Pakenham Campus Pakenham', 10)
    ,@('Covenant College Bell Post Hill', 24)
    ,@('Devon Meadows Primary School Devon Meadows', 12)
    ,@('Domestic Freight Vessel Siem Aquamarine Geelong', 10)
    ,@('Flemington Racecourse Flemington', 12)
    ,@('Gilly''s Early Learning Centre Balaclava', 11)
    ,@('Gladstone Park Secondary College 29 Oct Gladstone Park', 11)
    ,@('Guardian Childcare & Education Moorabbin', 12)
    ,@('Hamlyn Banks Primary School Hamlyn Heights', 10)
    ,@('Hazel Glen College Doreen', 14)
    ,@('Hazelwood North Primary School Hazelwood North', 27)
    ,@('Industries Pty Ltd Factory Thomastown', 17)
    ,@('Islamic College of Melbourne Tarneit Oct Nov', 26)
    ,@('Lyndhurst Primary School Lyndhurst', 12)
    ,@('Master Poultry Group West Footscray', 14)
    ,@('Minaret College Officer Campus Officer', 12)
    ,@('Morwell Park Primary School Morwell', 55)
    ,@('Narre Warren South P-12 College Narre Warren', 12)
    ,@('Nido Early School Woodend', 11)
    ,@('Northern Bay College Wexford Campus Corio', 27)
    ,@('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 14)
    ,@('Northern Health The Northern Hospital Epping', 15)
    ,@('Oakleigh South Primary School Oakleigh South', 17)
    ,@('Pentland Primary School Darley', 11)
    ,@('Rutherglen Motor Inn and Walkabout Motel Rutherglen', 22)
    ,@('Sirius College Ibrahim Dellal Campus Sunshine', 11)
    ,@('Sirius College Shepparton Campus Shepparton', 19)
    ,@('Smartie Pants Early Learning and Development Diamond Creek', 20)
    ,@('Social Gathering Woodvale 30 Oct', 10)
    ,@('Society Restaurant Melbourne', 24)
    ,@('St Ambrose Parish Primary School Woodend', 10)
    ,@('St Brendans Primary School Shepparton', 11)
    ,@('St Clare''s Primary School Officer', 12)
    ,@('St Georges Road Primary School Shepparton', 14)
    ,@('St Joseph''s School Quarry Hill', 33)
    ,@('St Louis de Montfort''s School Aspendale', 14)
    ,@('St Paul''s Primary School Sunshine West', 15)
    ,@('St Vincents Hospital Melbourne Emergency Department Fitzroy', 14)
    ,@('Stevensville Primary School St Albans', 10)
    ,@('Stockdale Road Primary School Traralgon', 33)
    ,@('Story House Early Learning Epping October', 11)
    ,@('Sunbury Primary School Sunbury', 11)
    ,@('TUROSI PTY LTD Thomastown', 10)
    ,@('Templestowe Park Primary School Templestowe', 31)
    ,@('The Lake Primary School Cabarita', 23)
    ,@('The Royal Children''s Hospital Melbourne Emergency Department Parkville', 15)
    ,@('Top Yard Rooftop Melbourne', 14)
    ,@('Tucker Road Bentleigh Primary School Bentleigh', 10)
    ,@('Warragul Regional College Warragul', 18)
    ,@('Werribee Mercy Hospital Emergency Department', 19)
    ,@('Western Health Sunshine Hospital Emergency Department St Albans', 13)
    ,@('Wodonga Primary School Wodonga', 14)
    ,@('Wodonga South Primary School Wodonga', 11)
    ,@('Wyndham Christian College Wyndham Vale', 12)
    ,@('Yeshivah College St Kilda East', 20)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}